# Update the "From" value of rule R30 (row 10) on the Rules sheet
# from 18 to 1, keeping its existing number style (s="20") intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
